# Auto update Excel log

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $text, $styleSourceCell) {
    # Force the cell to be treated as text so date-looking strings
    # (e.g. "2026-02-01") are not auto-converted into date serial numbers.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    # Re-apply the plain/default style (copied from an existing unstyled
    # cell) so no extra formatting is left behind on the cell.
    $cell.Style = $styleSourceCell.Style
}

# --- Proximity sheet: append rows 26-30 ---
$proximity = $wb.Worksheets.Item("Proximity")
$proxStyleSource = $proximity.Cells.Item(2, 4)

$proximityRows = @(
    @("2026-02-01", "15:10:45", "15:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "15:10:58", "15:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "15:11:08", "15:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "15:11:31", "15:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "15:11:31", "15:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
)

$startRow = 26
for ($i = 0; $i -lt $proximityRows.Count; $i++) {
    $row = $startRow + $i
    $data = $proximityRows[$i]
    Set-TextCell $proximity $row 1 $data[0] $proxStyleSource
    $proximity.Cells.Item($row, 2).Value = $data[1]
    $proximity.Cells.Item($row, 3).Value = $data[2]
    $proximity.Cells.Item($row, 4).Value = $data[3]
    $proximity.Cells.Item($row, 5).Value = $data[4]
    $proximity.Cells.Item($row, 6).Value = $data[5]
}

# --- Camera sheet: append rows 2-4 ---
$camera = $wb.Worksheets.Item("Camera")
$camStyleSource = $camera.Cells.Item(1, 4)

$cameraRows = @(
    @("2026-02-01", "15:10:57", "15:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "15:11:30", "15:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "15:11:33", "15:00", "Living Room Main Door", "Image Captured", "Active")
)

$startRow = 2
for ($i = 0; $i -lt $cameraRows.Count; $i++) {
    $row = $startRow + $i
    $data = $cameraRows[$i]
    Set-TextCell $camera $row 1 $data[0] $camStyleSource
    $camera.Cells.Item($row, 2).Value = $data[1]
    $camera.Cells.Item($row, 3).Value = $data[2]
    $camera.Cells.Item($row, 4).Value = $data[3]
    $camera.Cells.Item($row, 5).Value = $data[4]
    $camera.Cells.Item($row, 6).Value = $data[5]
}
